$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: row 3 is the 68d3c360-... file. Mark it handed back.
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# zh-cn sheet: row 3 is the 68d3c360-... file (Status + Latest Handback DateTime)
$zhcn.Range("B3").Value = $newStatus
$zhcn.Range("G3").Value = "2016-03-10 04:58:29"

# de-de sheet: row 3 is the 68d3c360-... file (Status + Latest Handback DateTime)
$dede.Range("B3").Value = $newStatus
$dede.Range("G3").Value = "2016-03-10 04:58:36"
